# Adds a new attendance column (J) to the "asistencia" sheet with a new
# date/header label "9/)" and per-student P/A values, a COUNTIF total in
# J28, then updates the active-sheet/selection state so "notas" becomes
# the selected sheet (matching the author's final on-screen state).

$wb = $excel.ActiveWorkbook

# --- asistencia: new column J -------------------------------------------------
$ws1 = $wb.Worksheets.Item("asistencia")

$ws1.Range("J1").Value = "9/)"

for ($r = 2; $r -le 27; $r++) {
    if ($r -eq 7) {
        $ws1.Cells.Item($r, 10).Value = "A"
    } else {
        $ws1.Cells.Item($r, 10).Value = "P"
    }
}

$ws1.Range("J28").Formula = "=COUNTIF(J2:J27,""P"")"

# Move asistencia's selection before switching sheets
$ws1.Range("A1").Select()
$ws1.Range("J7").Select()

# --- switch active sheet to notas, update its selection -----------------------
$ws2 = $wb.Worksheets.Item("notas")
$ws2.Activate()
$ws2.Range("N2").Select()
